# LambdaMART.xlsx edit:
#  - Add three new metric columns to the results table:
#      * "ERR@10 on train data:"       (inserted after "NDCG@10 on training data:")
#      * "ERR@10 on validation data:"  (inserted after "NDCG@10 on validation data:")
#      * "NDCG@10 on test data"        (inserted before "ERR@10 on test data:")
#  - Existing "NDCG@10 ..." / final "ERR@10 on test data:" value columns get a
#    bold look; all data cells become centered.
#  - Last header cell ("ERR@10 on test data:") gets its own bold header look.
#  - Selection left on C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Grow the table (and the backing worksheet range) from 4 to 7 columns.
#    Resize keeps the table's identity / style / per-column dxf wiring,
#    unlike delete+recreate.
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G6"))

# ---------------------------------------------------------------------------
# 2. Normalize formatting across the whole (now 7-wide) table first, by
#    cloning the existing header look (A1) across row 1, and the existing
#    data look (A2) across the data rows. This gives every new cell the
#    same base font/fill the original columns already had, before any
#    column-specific emphasis (bold / center) is layered on top.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2").Copy()
$ws.Range("A2:G6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Final column layout (A..G):
#      A Fold
#      B NDCG@10 on training data:
#      C ERR@10 on train data:        <- new
#      D NDCG@10 on validation data:
#      E ERR@10 on validation data:   <- new
#      F NDCG@10 on test data         <- new
#      G ERR@10 on test data:
# ---------------------------------------------------------------------------

$headers = @{
    "A1" = "Fold"
    "B1" = "NDCG@10 on training data:"
    "C1" = "ERR@10 on train data:"
    "D1" = "NDCG@10 on validation data:"
    "E1" = "ERR@10 on validation data:"
    "F1" = "NDCG@10 on test data"
    "G1" = "ERR@10 on test data:"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

$data = @{
    2 = @(1, 0.4633, 0.3561, 0.4569, 0.352,  0.4452, 0.3449)
    3 = @(2, 0.4762, 0.3674, 0.4466, 0.3438, 0.4461, 0.3371)
    4 = @(3, 0.479,  0.372,  0.4474, 0.3411, 0.4456, 0.334)
    5 = @(4, 0.4745, 0.3708, 0.4445, 0.3341, 0.4524, 0.3444)
    6 = @(5, 0.4728, 0.3662, 0.454,  0.3446, 0.4582, 0.3519)
}
$cols = @("A", "B", "C", "D", "E", "F", "G")
foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 4. Formatting touch-ups.
# ---------------------------------------------------------------------------

# All data cells (rows 2-6) are centered, horizontally and vertically.
$dataRange = $ws.Range("A2:G6")
$dataRange.HorizontalAlignment = -4108   # xlCenter
$dataRange.VerticalAlignment = -4108     # xlCenter

# The "headline" value columns (NDCG train / NDCG valid / ERR test) are bold.
$boldCols = @("B", "D", "G")
foreach ($c in $boldCols) {
    $ws.Range($c + "2:" + $c + "6").Font.Bold = $true
}

# The trailing header cell ("ERR@10 on test data:") gets a bold version of
# the default (dark-fill / black-font) header look, setting it apart from
# the rest of the header row.
$lastHeader = $ws.Range("G1")
$lastHeader.Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. Column widths (best-fit-like, matching the new layout).
# ---------------------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 22.5
$ws.Columns("D:D").ColumnWidth = 28.6640625
$ws.Columns("E:E").ColumnWidth = 28
$ws.Columns("F:F").ColumnWidth = 23
$ws.Columns("G:G").ColumnWidth = 21.83203125

# ---------------------------------------------------------------------------
# 6. Selection, as recorded in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("C3").Select()
